$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'69.069.57"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.51%  "

# Row 3
$ws.Range("D3").Value = "'3.754.51"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.03%  "

# Row 4
$ws.Range("D4").Value = "'0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.07%  "

# Row 5
$ws.Range("D5").Value = "'602.88"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.09%  "

# Row 6
$ws.Range("D6").Value = "'166.18"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.82%  "

# Row 7
$ws.Range("D7").Value = "'3.749.83"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.11%  "

# Row 9
$ws.Range("E9").Value = "  +0.52%  "

# Row 10
$ws.Range("E10").Value = "  +5.16%  "

# Row 11
$ws.Range("D11").Value = "'6.36"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.57%  "

# Row 12
$ws.Range("E12").Value = "  -0.85%  "

# Row 13
$ws.Range("D13").Value = "'37.72"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.62%  "

# Row 14
$ws.Range("D14").Value = "'0.0000247"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.57%  "

# Row 15
$ws.Range("D15").Value = "'4.380.51"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.08%  "

# Row 16
$ws.Range("D16").Value = "'3.754.16"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.01%  "

# Row 17
$ws.Range("D17").Value = "'69.027.02"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.40%  "

# Row 18
$ws.Range("D18").Value = "'7.43"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.55%  "

# Row 19
$ws.Range("D19").Value = "'17.83"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.96%  "

# Row 21
$ws.Range("D21").Value = "'11.14"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.74%  "

# Row 22
$ws.Range("D22").Value = "'490.60"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.32%  "

# Row 23
$ws.Range("D23").Value = "'0.724"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.72%  "

# Row 25
$ws.Range("E25").Value = "  -1.35%  "

# Row 26
$ws.Range("E26").Value = "  -2.59%  "

# Row 27
$ws.Range("D27").Value = "'12.29"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.47%  "

# Row 28
$ws.Range("D28").Value = "'10.04"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.91%  "

# Row 30
$ws.Range("E30").Value = "  -0.75%  "

# Row 31
$ws.Range("D31").Value = "'8.13"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.02%  "

# Row 32
$ws.Range("D32").Value = "'2.43"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.96%  "

# Row 33
$ws.Range("D33").Value = "'31.60"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.59%  "

# Row 34
$ws.Range("D34").Value = "'3.894.72"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.21%  "

# Row 35
$ws.Range("D35").Value = "'3.691.89"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.07%  "

# Row 36
$ws.Range("E36").Value = "  -0.56%  "

# Row 37
$ws.Range("E37").Value = "  +5.79%  "

# Row 38
$ws.Range("B38").Value = "Mantle"
$ws.Range("C38").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D38").Value = "'1.01"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.59%  "

# Row 39
$ws.Range("B39").Value = "Filecoin"
$ws.Range("C39").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D39").Value = "'5.93"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.39%  "

# Row 40
$ws.Range("E40").Value = "  -0.07%  "

# Row 41
$ws.Range("E41").Value = "  +8.80%  "

# Row 42
$ws.Range("D42").Value = "'0.325"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.17%  "

# Row 43
$ws.Range("D43").Value = "'48.51"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.74%  "

# Row 44
$ws.Range("B44").Value = "Bittensor"
$ws.Range("C44").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D44").Value = "'425.79"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -4.29%  "

# Row 45
$ws.Range("B45").Value = "Stacks"
$ws.Range("C45").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D45").Value = "'1.99"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.46%  "

# Row 46
$ws.Range("B46").Value = "USDe"
$ws.Range("C46").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D46").Value = "'1.00"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.02%  "

# Row 47
$ws.Range("B47").Value = "Cosmos"
$ws.Range("C47").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D47").Value = "'8.42"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.97%  "

# Row 48
$ws.Range("B48").Value = "Arweave"
$ws.Range("C48").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D48").Value = "'40.10"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.18%  "

# Row 49
$ws.Range("B49").Value = "Monero"
$ws.Range("C49").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D49").Value = "'142.34"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.11%  "

# Row 50
$ws.Range("E50").Value = "  +7.41%  "

# Row 51
$ws.Range("D51").Value = "'2.787.46"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.11%  "
